# Refresh the cryptos.xlsx price/volume snapshot (scheduled GitHub Actions run).
# Row 2-51 hold the coin table; columns are B=Coin, C=Link, D=Price, E=Volume(1h).
# Rows 48/49 additionally swap their VeChain / EnergySwap content.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '57.783.73'  # D2
$ws.Cells.Item(2, 5).Value = '  -1.43%  '  # E2

$ws.Cells.Item(3, 4).Value = '2.446.83'  # D3
$ws.Cells.Item(3, 5).Value = '  -2.93%  '  # E3

$ws.Cells.Item(4, 5).Value = '  +0.47%  '  # E4

$ws.Cells.Item(5, 4).Value = '''523.57'  # D5
$ws.Cells.Item(5, 5).Value = '  +0.38%  '  # E5

$ws.Cells.Item(6, 4).Value = '''129.76'  # D6
$ws.Cells.Item(6, 5).Value = '  -2.25%  '  # E6

$ws.Cells.Item(7, 4).Value = '''1.00'  # D7
$ws.Cells.Item(7, 5).Value = '  +0.62%  '  # E7

$ws.Cells.Item(8, 5).Value = '  +0.70%  '  # E8

$ws.Cells.Item(9, 4).Value = '''0.0976'  # D9
$ws.Cells.Item(9, 5).Value = '  +0.01%  '  # E9

$ws.Cells.Item(10, 5).Value = '  -3.51%  '  # E10

$ws.Cells.Item(11, 4).Value = '''4.93'  # D11
$ws.Cells.Item(11, 5).Value = '  -4.19%  '  # E11

$ws.Cells.Item(12, 4).Value = '''0.321'  # D12
$ws.Cells.Item(12, 5).Value = '  -3.39%  '  # E12

$ws.Cells.Item(13, 4).Value = '2.880.66'  # D13
$ws.Cells.Item(13, 5).Value = '  -1.81%  '  # E13

$ws.Cells.Item(14, 4).Value = '57.706.71'  # D14
$ws.Cells.Item(14, 5).Value = '  -1.12%  '  # E14

$ws.Cells.Item(15, 4).Value = '''21.69'  # D15
$ws.Cells.Item(15, 5).Value = '  -1.94%  '  # E15

$ws.Cells.Item(16, 4).Value = '''0.0000132'  # D16
$ws.Cells.Item(16, 5).Value = '  -2.03%  '  # E16

$ws.Cells.Item(17, 4).Value = '2.449.87'  # D17
$ws.Cells.Item(17, 5).Value = '  -2.08%  '  # E17

$ws.Cells.Item(18, 4).Value = '''10.35'  # D18
$ws.Cells.Item(18, 5).Value = '  -3.04%  '  # E18

$ws.Cells.Item(19, 4).Value = '''4.13'  # D19
$ws.Cells.Item(19, 5).Value = '  -0.77%  '  # E19

$ws.Cells.Item(20, 4).Value = '''314.82'  # D20
$ws.Cells.Item(20, 5).Value = '  -2.32%  '  # E20

$ws.Cells.Item(21, 5).Value = '  +0.20%  '  # E21

$ws.Cells.Item(22, 5).Value = '  +0.42%  '  # E22

$ws.Cells.Item(23, 4).Value = '''64.91'  # D23
$ws.Cells.Item(23, 5).Value = '  +0.84%  '  # E23

$ws.Cells.Item(24, 5).Value = '  +1.63%  '  # E24

$ws.Cells.Item(25, 5).Value = '  +1.34%  '  # E25

$ws.Cells.Item(26, 5).Value = '  -2.62%  '  # E26

$ws.Cells.Item(27, 4).Value = '''7.21'  # D27
$ws.Cells.Item(27, 5).Value = '  -2.29%  '  # E27

$ws.Cells.Item(28, 4).Value = '''172.16'  # D28
$ws.Cells.Item(28, 5).Value = '  +2.56%  '  # E28

$ws.Cells.Item(29, 4).Value = '0.0₃0733'  # D29
$ws.Cells.Item(29, 5).Value = '  -2.96%  '  # E29

$ws.Cells.Item(30, 5).Value = '  -1.36%  '  # E30

$ws.Cells.Item(31, 5).Value = '  -3.97%  '  # E31

$ws.Cells.Item(32, 4).Value = '''6.09'  # D32
$ws.Cells.Item(32, 5).Value = '  -2.68%  '  # E32

$ws.Cells.Item(33, 5).Value = '  +0.13%  '  # E33

$ws.Cells.Item(34, 5).Value = '  +0.44%  '  # E34

$ws.Cells.Item(35, 4).Value = '''17.79'  # D35
$ws.Cells.Item(35, 5).Value = '  -1.85%  '  # E35

$ws.Cells.Item(36, 5).Value = '  -6.83%  '  # E36

$ws.Cells.Item(37, 4).Value = '''3.80'  # D37
$ws.Cells.Item(37, 5).Value = '  -4.34%  '  # E37

$ws.Cells.Item(38, 4).Value = '''36.27'  # D38
$ws.Cells.Item(38, 5).Value = '  +0.90%  '  # E38

$ws.Cells.Item(39, 4).Value = '''1.46'  # D39
$ws.Cells.Item(39, 5).Value = '  -0.65%  '  # E39

$ws.Cells.Item(40, 4).Value = '''0.789'  # D40
$ws.Cells.Item(40, 5).Value = '  +1.26%  '  # E40

$ws.Cells.Item(41, 4).Value = '''3.40'  # D41
$ws.Cells.Item(41, 5).Value = '  -2.79%  '  # E41

$ws.Cells.Item(42, 4).Value = '''264.01'  # D42
$ws.Cells.Item(42, 5).Value = '  -5.04%  '  # E42

$ws.Cells.Item(43, 4).Value = '''0.582'  # D43
$ws.Cells.Item(43, 5).Value = '  -2.70%  '  # E43

$ws.Cells.Item(44, 4).Value = '''4.79'  # D44
$ws.Cells.Item(44, 5).Value = '  -5.82%  '  # E44

$ws.Cells.Item(45, 4).Value = '''124.11'  # D45
$ws.Cells.Item(45, 5).Value = '  +0.50%  '  # E45

$ws.Cells.Item(46, 4).Value = '''0.0927'  # D46
$ws.Cells.Item(46, 5).Value = '  +0.93%  '  # E46

$ws.Cells.Item(47, 4).Value = '''0.0492'  # D47
$ws.Cells.Item(47, 5).Value = '  -1.78%  '  # E47

$ws.Cells.Item(48, 2).Value = 'EnergySwap'  # B48
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'  # C48
$ws.Cells.Item(48, 4).Value = '''17.03'  # D48
$ws.Cells.Item(48, 5).Value = '  -4.55%  '  # E48

$ws.Cells.Item(49, 2).Value = 'VeChain'  # B49
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'  # C49
$ws.Cells.Item(49, 4).Value = '''0.0210'  # D49
$ws.Cells.Item(49, 5).Value = '  -1.89%  '  # E49

$ws.Cells.Item(50, 4).Value = '''16.28'  # D50
$ws.Cells.Item(50, 5).Value = '  -3.94%  '  # E50

$ws.Cells.Item(51, 4).Value = '1.706.51'  # D51
$ws.Cells.Item(51, 5).Value = '  -2.08%  '  # E51
